# Fruta / hortaliza, semanal
#
# Inserts 4 new weekly price rows (Durazno / Femacal de La Calera) right
# before the existing row 430, pushing the rest of the table down by 4
# rows (old rows 430-496 become new rows 434-500). The sheet's used
# range grows from A1:T496 to A1:T500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 430, shifting everything below down.
$ws.Rows("430:433").Insert()

function Set-PriceRow($r, $mercado, $region, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $origen, $precioKg, $kgUnidad) {
    $ws.Cells.Item($r, 1).Value = 3
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = 5
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100103
    $ws.Cells.Item($r, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($r, 9).Value = 100103004
    $ws.Cells.Item($r, 10).Value = "Durazno"
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-PriceRow 430 "Femacal de La Calera" "Coquimbo" 44504 "Early Majestic" "Segunda" 70 12000 12000 12000 "`$/bandeja 10 kilos granel" "Provincia de San Felipe de Aconcagua" 1200 10

Set-PriceRow 431 "Femacal de La Calera" "Coquimbo" 44504 "Florida King" "Primera" 56 13000 13000 13000 "`$/bandeja 10 kilos empedrada" "Provincia de San Felipe de Aconcagua" 1300 10

Set-PriceRow 432 "Femacal de La Calera" "Coquimbo" 44504 "Florida King" "Primera" 68 17000 17000 17000 "`$/caja 15 kilos granel" "Provincia de Limarí" 1133 15

Set-PriceRow 433 "Femacal de La Calera" "Coquimbo" 44504 "Florida King" "Segunda" 67 14000 14000 14000 "`$/caja 15 kilos granel" "Provincia de Limarí" 933 15

# Make sure the date column keeps its date number-format on the new rows
# (Insert() above already copies it from the neighbouring row, this is
# just a safety net).
$ws.Range("D430:D433").NumberFormat = "YYYY-MM-DD HH:MM:SS"
